$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D and E
$ws.Range("D2:D3").Value = 0.133
$ws.Range("E2:E3").Value = 0.146

# Columns I, J
$ws.Range("I2:I3").Value = 0
$ws.Range("J2:J3").Value = 0

# Column K
$ws.Range("K2:K3").Value = 114.9

# Column L
$ws.Range("L2:L3").Value = 0.4368821292775666

# Column M
$ws.Range("M2:M3").Value = 0.606

# Column N
$ws.Range("N2:N3").Value = 0.0006530172413793103

# Column O
$ws.Range("O2:O3").Value = 0.005274151436031331

# Column P
$ws.Range("P2:P3").Value = 0.606

# Column Q
$ws.Range("Q2:Q3").Value = 0.0006530172413793103

# Column R
$ws.Range("R2:R3").Value = 0.005274151436031331

# Column U
$ws.Range("U2:U3").Value = 442.5

# Column V
$ws.Range("V2:V3").Value = 0.4768318965517241

# Column W
$ws.Range("W2:W3").Value = 0.136250444681608

# Column X
$ws.Range("X2:X3").Value = 0.1111128840395408

# Column Y
$ws.Range("Y2:Y3").Value = 0.02513756064206718

# Column Z
$ws.Range("Z2:Z3").Value = 0.1173635592842162

# Column AA
$ws.Range("AA2:AA3").Value = 0

# Column AB
$ws.Range("AB2:AB3").Value = 0.06599220328306395

# Column AC
$ws.Range("AC2:AC3").Value = -0.06599220328306395

# Column AD
$ws.Range("AD2:AD3").Value = 1947.9

# Column AE
$ws.Range("AE2:AE3").Value = 0

# Column AF
$ws.Range("AF2:AF3").Value = 1947.9

# Column AG
$ws.Range("AG2:AG3").Value = 1505.4

# Column AH
$ws.Range("AH2:AH3").Value = 0.6773184046733197

# Column AI
$ws.Range("AI2:AI3").Value = 0.6937954124519161

# Column AJ
$ws.Range("AJ2:AJ3").Value = 0.6186405851894469

# Column AK
$ws.Range("AK2:AK3").Value = 0.6365058559891759

# Columns AN and AP are removed entirely in rows 2 and 3
$ws.Range("AN2:AN3").ClearContents()
$ws.Range("AP2:AP3").ClearContents()
